$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" numeric-looking string need to be forced to Text
# (matching the source data which stores these as literal strings), otherwise Excel
# auto-converts them to numbers on assignment.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2
$ws.Range("D2").Value = "68.255.30"
$ws.Range("E2").Value = "  +2.25%  "

# Row 3
$ws.Range("D3").Value = "3.132.59"
$ws.Range("E3").Value = "  +2.00%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
Set-TextValue "D5" "576.99"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
Set-TextValue "D6" "180.60"
$ws.Range("E6").Value = "  +6.06%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "3.131.15"
$ws.Range("E8").Value = "  +2.01%  "

# Row 9
$ws.Range("E9").Value = "  +1.73%  "

# Row 10
Set-TextValue "D10" "6.53"
$ws.Range("E10").Value = "  +2.47%  "

# Row 11
$ws.Range("E11").Value = "  +2.05%  "

# Row 12
Set-TextValue "D12" "0.470"
$ws.Range("E12").Value = "  +0.70%  "

# Row 13
$ws.Range("E13").Value = "  +1.24%  "

# Row 14
Set-TextValue "D14" "36.73"
$ws.Range("E14").Value = "  +2.67%  "

# Row 16
$ws.Range("D16").Value = "68.256.89"
$ws.Range("E16").Value = "  +2.24%  "

# Row 17
$ws.Range("D17").Value = "3.658.67"
$ws.Range("E17").Value = "  +2.02%  "

# Row 18
Set-TextValue "D18" "7.11"
$ws.Range("E18").Value = "  +1.75%  "

# Row 19
$ws.Range("D19").Value = "3.134.75"
$ws.Range("E19").Value = "  +2.20%  "

# Row 20
Set-TextValue "D20" "16.60"
$ws.Range("E20").Value = "  -1.96%  "

# Row 21
Set-TextValue "D21" "487.53"
$ws.Range("E21").Value = "  -0.65%  "

# Row 22
Set-TextValue "D22" "0.697"
$ws.Range("E22").Value = "  +1.45%  "

# Row 23
Set-TextValue "D23" "7.79"
$ws.Range("E23").Value = "  +1.11%  "

# Row 24
Set-TextValue "D24" "83.88"
$ws.Range("E24").Value = "  +1.24%  "

# Row 25
Set-TextValue "D25" "12.97"
$ws.Range("E25").Value = "  +2.51%  "

# Row 26
Set-TextValue "D26" "2.33"
$ws.Range("E26").Value = "  +5.88%  "

# Row 27
$ws.Range("E27").Value = "  +4.61%  "

# Row 28
$ws.Range("E28").Value = "  +0.01%  "

# Row 29
Set-TextValue "D29" "8.12"
$ws.Range("E29").Value = "  +4.29%  "

# Row 30
$ws.Range("E30").Value = "  +4.22%  "

# Row 31
$ws.Range("E31").Value = "  +1.09%  "

# Row 32
Set-TextValue "D32" "28.27"
$ws.Range("E32").Value = "  +2.64%  "

# Row 33
$ws.Range("E33").Value = "  +1.14%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0952"
$ws.Range("E34").Value = "  +4.37%  "

# Row 35
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
Set-TextValue "D36" "48.95"
$ws.Range("E36").Value = "  +4.23%  "

# Row 37
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D37" "0.958"
$ws.Range("E37").Value = "  +1.00%  "

# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D38" "5.65"
$ws.Range("E38").Value = "  +1.21%  "

# Row 39
Set-TextValue "D39" "0.324"
$ws.Range("E39").Value = "  +8.53%  "

# Row 40
Set-TextValue "D40" "2.05"
$ws.Range("E40").Value = "  +4.60%  "

# Row 41
Set-TextValue "D41" "49.11"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("E42").Value = "  +1.22%  "

# Row 43
Set-TextValue "D43" "8.42"
$ws.Range("E43").Value = "  +1.31%  "

# Row 44
$ws.Range("E44").Value = "  +8.48%  "

# Row 45
Set-TextValue "D45" "395.85"
$ws.Range("E45").Value = "  +8.24%  "

# Row 46
$ws.Range("D46").Value = "2.785.85"
$ws.Range("E46").Value = "  +1.20%  "

# Row 47
Set-TextValue "D47" "27.06"
$ws.Range("E47").Value = "  +9.74%  "

# Row 48
Set-TextValue "D48" "0.0349"
$ws.Range("E48").Value = "  +1.03%  "

# Row 49
Set-TextValue "D49" "135.65"
$ws.Range("E49").Value = "  +0.24%  "

# Row 50
$ws.Range("E50").Value = "  +0.06%  "

# Row 51
Set-TextValue "D51" "2.34"
$ws.Range("E51").Value = "  +8.72%  "
